# Update "想去人数" (want-to-go count) figures on the 苏州-漫展信息 workbook
# to match newly scraped totals (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 812
$wsExhibit.Range("F7").Value = 242
$wsExhibit.Range("F14").Value = 12858
$wsExhibit.Range("F16").Value = 5279
$wsExhibit.Range("F17").Value = 5535

# --- Sheet "演出" (performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 115

# --- Sheet "全部类型" (all types, combined listing) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 812
$wsAll.Range("F7").Value = 242
$wsAll.Range("F14").Value = 12858
$wsAll.Range("F15").Value = 115
$wsAll.Range("F18").Value = 5279
$wsAll.Range("F19").Value = 5535
